$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()

# Update the "Do HW2" task hours for G10:G12 (member3 sprint day 3 column),
# merging in an extra hour: the values become hard literals (11) instead of
# the previous shared fill-formula ("=F10" etc.), matching the edit where the
# author typed a new value directly into the merged cells.
$ws.Range("G10:G12").Value = 11

# Reflect the cell(s) the author was working in when the change was made.
$ws.Range("G10:G12").Select()
$ws.Application.ActiveWindow.RangeSelection.Item(1).Activate()

$excel.ActiveWindow.Height = 7680
